$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 2159.92  # H51: 2176 -> 2159.92
$ws.Cells.Item(51, 9).Value = 1926.6666  # I51: 1931.25 -> 1926.6666
$ws.Cells.Item(51, 10).Value = 2509.8  # J51: 2611.111 -> 2509.8
$ws.Cells.Item(51, 11).Value = 1926.6666  # K51: 1931.25 -> 1926.6666
$ws.Cells.Item(51, 12).Value = 2509.8  # L51: 2611.111 -> 2509.8
$ws.Cells.Item(51, 13).Value = -1442.6666  # M51: -1447.25 -> -1442.6666
$ws.Cells.Item(51, 14).Value = -3477.8  # N51: -3579.111 -> -3477.8

$ws.Cells.Item(86, 8).Value = 5551.077  # H86: 5753.16 -> 5551.077
$ws.Cells.Item(86, 9).Value = 13439.25  # I86: 15287.857 -> 13439.25
$ws.Cells.Item(86, 11).Value = 13439.25  # K86: 15287.857 -> 13439.25
$ws.Cells.Item(86, 13).Value = -12316.25  # M86: -14164.857 -> -12316.25

$ws.Cells.Item(89, 8).Value = 5551.077  # H89: 5753.16 -> 5551.077
$ws.Cells.Item(89, 9).Value = 13439.25  # I89: 15287.857 -> 13439.25
$ws.Cells.Item(89, 11).Value = 67196.25  # K89: 76439.285 -> 67196.25
$ws.Cells.Item(89, 13).Value = -61580.25  # M89: -70823.285 -> -61580.25

$ws.Cells.Item(106, 8).Value = 856.8182  # H106: 868.75 -> 856.8182
$ws.Cells.Item(106, 9).Value = 856.8182  # I106: 868.75 -> 856.8182
$ws.Cells.Item(106, 11).Value = 856.8182  # K106: 868.75 -> 856.8182
$ws.Cells.Item(106, 13).Value = -225.8182  # M106: -237.75 -> -225.8182

$ws.Cells.Item(125, 8).Value = 1375  # H125: 2766.6667 -> 1375
$ws.Cells.Item(125, 9).Value = 880.5714  # I125: 1300 -> 880.5714
$ws.Cells.Item(125, 10).Value = 2067.2  # J125: 3500 -> 2067.2
$ws.Cells.Item(125, 11).Value = 7925.1426  # K125: 11700 -> 7925.1426
$ws.Cells.Item(125, 12).Value = 18604.8  # L125: 31500 -> 18604.8
$ws.Cells.Item(125, 13).Value = -5465.1426  # M125: -9240 -> -5465.1426
$ws.Cells.Item(125, 14).Value = -23524.8  # N125: -36420 -> -23524.8

$ws.Cells.Item(129, 8).Value = 4099910.8  # H129: 4099911 -> 4099910.8
$ws.Cells.Item(129, 9).Value = 27778762  # I129: 27778764 -> 27778762
$ws.Cells.Item(129, 11).Value = 83336286  # K129: 83336292 -> 83336286
$ws.Cells.Item(129, 13).Value = -83331286  # M129: -83331292 -> -83331286

$ws.Cells.Item(132, 8).Value = 18184822  # H132: 8336278 -> 18184822
$ws.Cells.Item(132, 9).Value = 22223672  # I132: 9093213 -> 22223672
$ws.Cells.Item(132, 11).Value = 66671016  # K132: 27279639 -> 66671016
$ws.Cells.Item(132, 13).Value = -66668486  # M132: -27277109 -> -66668486

$ws.Cells.Item(137, 8).Value = 5268331.5  # H137: 4549999.5 -> 5268331.5
$ws.Cells.Item(137, 9).Value = 9097391  # I137: 6671666 -> 9097391
$ws.Cells.Item(137, 10).Value = 3375.375  # J137: 3571.8572 -> 3375.375
$ws.Cells.Item(137, 11).Value = 27292173  # K137: 20014998 -> 27292173
$ws.Cells.Item(137, 12).Value = 10126.125  # L137: 10715.5716 -> 10126.125
$ws.Cells.Item(137, 13).Value = -27289623  # M137: -20012448 -> -27289623
$ws.Cells.Item(137, 14).Value = -15226.125  # N137: -15815.5716 -> -15226.125

$ws.Cells.Item(138, 8).Value = 3290.3677  # H138: 3360.7693 -> 3290.3677
$ws.Cells.Item(138, 9).Value = 1465.7576  # I138: 1489.7812 -> 1465.7576
$ws.Cells.Item(138, 10).Value = 5010.7144  # J138: 5175.0605 -> 5010.7144
$ws.Cells.Item(138, 11).Value = 4397.2728  # K138: 4469.3436 -> 4397.2728
$ws.Cells.Item(138, 12).Value = 15032.1432  # L138: 15525.1815 -> 15032.1432
$ws.Cells.Item(138, 13).Value = 742.7272000000003  # M138: 670.6563999999998 -> 742.7272000000003
$ws.Cells.Item(138, 14).Value = -25312.1432  # N138: -25805.1815 -> -25312.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8146.447  # H32: 4764.9277 -> 8146.447
$ws.Cells.Item(32, 9).Value = 5702.1284  # I32: 3133.1868 -> 5702.1284
$ws.Cells.Item(32, 11).Value = 5702.1284  # K32: 3133.1868 -> 5702.1284
$ws.Cells.Item(32, 13).Value = -5415.1284  # M32: -2846.1868 -> -5415.1284

$ws.Cells.Item(88, 8).Value = 2451.2727  # H88: 2250 -> 2451.2727
$ws.Cells.Item(88, 10).Value = 2620.5  # J88: 2500 -> 2620.5
$ws.Cells.Item(88, 12).Value = 2620.5  # L88: 2500 -> 2620.5
$ws.Cells.Item(88, 14).Value = -3432.5  # N88: -3312 -> -3432.5

$ws.Cells.Item(91, 8).Value = 2451.2727  # H91: 2250 -> 2451.2727
$ws.Cells.Item(91, 10).Value = 2620.5  # J91: 2500 -> 2620.5
$ws.Cells.Item(91, 12).Value = 2620.5  # L91: 2500 -> 2620.5
$ws.Cells.Item(91, 14).Value = -5428.5  # N91: -5308 -> -5428.5

$ws.Cells.Item(97, 8).Value = 282.56  # H97: 282.4074 -> 282.56
$ws.Cells.Item(97, 9).Value = 289.73914  # I97: 292.70834 -> 289.73914
$ws.Cells.Item(97, 11).Value = 289.73914  # K97: 292.70834 -> 289.73914
$ws.Cells.Item(97, 13).Value = 206.26086  # M97: 203.29166 -> 206.26086

$ws.Cells.Item(122, 8).Value = 2708.6667  # H122: 2221.889 -> 2708.6667
$ws.Cells.Item(122, 9).Value = 1729.6666  # I122: 1733.5 -> 1729.6666
$ws.Cells.Item(122, 10).Value = 4666.6665  # J122: 2612.6 -> 4666.6665
$ws.Cells.Item(122, 11).Value = 5188.9998  # K122: 5200.5 -> 5188.9998
$ws.Cells.Item(122, 12).Value = 13999.9995  # L122: 7837.799999999999 -> 13999.9995
$ws.Cells.Item(122, 13).Value = -2738.9998  # M122: -2750.5 -> -2738.9998
$ws.Cells.Item(122, 14).Value = -18899.9995  # N122: -12737.8 -> -18899.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2020.9445  # H20: 1341.2258 -> 2020.9445
$ws.Cells.Item(20, 9).Value = 1161.1818  # I20: 821.8889 -> 1161.1818
$ws.Cells.Item(20, 10).Value = 3372  # J20: 2060.3076 -> 3372
$ws.Cells.Item(20, 11).Value = 1161.1818  # K20: 821.8889 -> 1161.1818
$ws.Cells.Item(20, 12).Value = 3372  # L20: 2060.3076 -> 3372
$ws.Cells.Item(20, 13).Value = -914.1818000000001  # M20: -574.8889 -> -914.1818000000001
$ws.Cells.Item(20, 14).Value = -3866  # N20: -2554.3076 -> -3866

$ws.Cells.Item(55, 8).Value = 0  # H55: 36000 -> 0
$ws.Cells.Item(55, 9).Value = 0  # I55: 28000 -> 0
$ws.Cells.Item(55, 10).Value = 0  # J55: 40000 -> 0
$ws.Cells.Item(55, 11).Value = 0  # K55: 28000 -> 0
$ws.Cells.Item(55, 12).ClearContents()  # L55: was 40000
$ws.Cells.Item(55, 13).ClearContents()  # M55: was -27727
$ws.Cells.Item(55, 14).Value = 0  # N55: -40546 -> 0

$ws.Cells.Item(86, 8).Value = 1941.5769  # H86: 2003.2 -> 1941.5769
$ws.Cells.Item(86, 9).Value = 1271.1666  # I86: 1298.6666 -> 1271.1666
$ws.Cells.Item(86, 10).Value = 3450  # J86: 3060 -> 3450
$ws.Cells.Item(86, 11).Value = 1271.1666  # K86: 1298.6666 -> 1271.1666
$ws.Cells.Item(86, 12).Value = 3450  # L86: 3060 -> 3450
$ws.Cells.Item(86, 13).Value = -148.1666  # M86: -175.6666 -> -148.1666
$ws.Cells.Item(86, 14).Value = -5696  # N86: -5306 -> -5696

$ws.Cells.Item(89, 8).Value = 1941.5769  # H89: 2003.2 -> 1941.5769
$ws.Cells.Item(89, 9).Value = 1271.1666  # I89: 1298.6666 -> 1271.1666
$ws.Cells.Item(89, 10).Value = 3450  # J89: 3060 -> 3450
$ws.Cells.Item(89, 11).Value = 6355.833000000001  # K89: 6493.333000000001 -> 6355.833000000001
$ws.Cells.Item(89, 12).Value = 17250  # L89: 15300 -> 17250
$ws.Cells.Item(89, 13).Value = -739.8330000000005  # M89: -877.3330000000005 -> -739.8330000000005
$ws.Cells.Item(89, 14).Value = -28482  # N89: -26532 -> -28482

$ws.Cells.Item(105, 8).Value = 1570.1034  # H105: 1681.9259 -> 1570.1034
$ws.Cells.Item(105, 9).Value = 1379.4445  # I105: 1531.25 -> 1379.4445
$ws.Cells.Item(105, 10).Value = 1882.091  # J105: 1901.091 -> 1882.091
$ws.Cells.Item(105, 11).Value = 1379.4445  # K105: 1531.25 -> 1379.4445
$ws.Cells.Item(105, 12).Value = 1882.091  # L105: 1901.091 -> 1882.091
$ws.Cells.Item(105, 13).Value = 367.5554999999999  # M105: 215.75 -> 367.5554999999999
$ws.Cells.Item(105, 14).Value = -5376.091  # N105: -5395.091 -> -5376.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3129997  # H31: 2862558.8 -> 3129997
$ws.Cells.Item(31, 9).Value = 5557500.5  # I31: 5557506 -> 5557500.5
$ws.Cells.Item(31, 10).Value = 8921.357  # J31: 9085.235000000001 -> 8921.357
$ws.Cells.Item(31, 11).Value = 5557500.5  # K31: 5557506 -> 5557500.5
$ws.Cells.Item(31, 12).Value = 8921.357  # L31: 9085.235000000001 -> 8921.357
$ws.Cells.Item(31, 13).Value = -5557205.5  # M31: -5557211 -> -5557205.5
$ws.Cells.Item(31, 14).Value = -9511.357  # N31: -9675.235000000001 -> -9511.357

$ws.Cells.Item(34, 8).Value = 3129997  # H34: 2862558.8 -> 3129997
$ws.Cells.Item(34, 9).Value = 5557500.5  # I34: 5557506 -> 5557500.5
$ws.Cells.Item(34, 10).Value = 8921.357  # J34: 9085.235000000001 -> 8921.357
$ws.Cells.Item(34, 11).Value = 5557500.5  # K34: 5557506 -> 5557500.5
$ws.Cells.Item(34, 12).Value = 8921.357  # L34: 9085.235000000001 -> 8921.357
$ws.Cells.Item(34, 13).Value = -5557298.5  # M34: -5557304 -> -5557298.5
$ws.Cells.Item(34, 14).Value = -9325.357  # N34: -9489.235000000001 -> -9325.357

$ws.Cells.Item(62, 8).Value = 4720.625  # H62: 4986.875 -> 4720.625
$ws.Cells.Item(62, 9).Value = 3573.3333  # I62: 4022.8572 -> 3573.3333
$ws.Cells.Item(62, 10).Value = 6195.7144  # J62: 5736.6665 -> 6195.7144
$ws.Cells.Item(62, 11).Value = 3573.3333  # K62: 4022.8572 -> 3573.3333
$ws.Cells.Item(62, 12).Value = 6195.7144  # L62: 5736.6665 -> 6195.7144
$ws.Cells.Item(62, 13).Value = -2949.3333  # M62: -3398.8572 -> -2949.3333
$ws.Cells.Item(62, 14).Value = -7443.7144  # N62: -6984.6665 -> -7443.7144

$ws.Cells.Item(65, 8).Value = 4720.625  # H65: 4986.875 -> 4720.625
$ws.Cells.Item(65, 9).Value = 3573.3333  # I65: 4022.8572 -> 3573.3333
$ws.Cells.Item(65, 10).Value = 6195.7144  # J65: 5736.6665 -> 6195.7144
$ws.Cells.Item(65, 11).Value = 17866.6665  # K65: 20114.286 -> 17866.6665
$ws.Cells.Item(65, 12).Value = 30978.572  # L65: 28683.3325 -> 30978.572
$ws.Cells.Item(65, 13).Value = -14746.6665  # M65: -16994.286 -> -14746.6665
$ws.Cells.Item(65, 14).Value = -37218.572  # N65: -34923.3325 -> -37218.572

$ws.Cells.Item(105, 8).Value = 1892.4375  # H105: 1986.0667 -> 1892.4375
$ws.Cells.Item(105, 9).Value = 1317.1818  # I105: 1343.6364 -> 1317.1818
$ws.Cells.Item(105, 10).Value = 3158  # J105: 3752.75 -> 3158
$ws.Cells.Item(105, 11).Value = 1317.1818  # K105: 1343.6364 -> 1317.1818
$ws.Cells.Item(105, 12).Value = 3158  # L105: 3752.75 -> 3158
$ws.Cells.Item(105, 13).Value = 429.8181999999999  # M105: 403.3635999999999 -> 429.8181999999999
$ws.Cells.Item(105, 14).Value = -6652  # N105: -7246.75 -> -6652

$ws.Cells.Item(122, 8).Value = 2034.3077  # H122: 2967.875 -> 2034.3077
$ws.Cells.Item(122, 9).Value = 1087.9  # I122: 2108.7778 -> 1087.9
$ws.Cells.Item(122, 10).Value = 5189  # J122: 4072.4285 -> 5189
$ws.Cells.Item(122, 11).Value = 3263.7  # K122: 6326.3334 -> 3263.7
$ws.Cells.Item(122, 12).Value = 15567  # L122: 12217.2855 -> 15567
$ws.Cells.Item(122, 13).Value = -813.7000000000003  # M122: -3876.3334 -> -813.7000000000003
$ws.Cells.Item(122, 14).Value = -20467  # N122: -17117.2855 -> -20467

$ws.Cells.Item(132, 8).Value = 3411.6897  # H132: 3285.516 -> 3411.6897
$ws.Cells.Item(132, 9).Value = 1935.0667  # I132: 1935.8667 -> 1935.0667
$ws.Cells.Item(132, 10).Value = 4993.7856  # J132: 4550.8125 -> 4993.7856
$ws.Cells.Item(132, 11).Value = 5805.2001  # K132: 5807.6001 -> 5805.2001
$ws.Cells.Item(132, 12).Value = 14981.3568  # L132: 13652.4375 -> 14981.3568
$ws.Cells.Item(132, 13).Value = -3275.2001  # M132: -3277.6001 -> -3275.2001
$ws.Cells.Item(132, 14).Value = -20041.3568  # N132: -18712.4375 -> -20041.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1489.5217  # H5: 1580.4286 -> 1489.5217
$ws.Cells.Item(5, 9).Value = 747.2941  # I5: 775.6 -> 747.2941
$ws.Cells.Item(5, 11).Value = 2241.8823  # K5: 2326.8 -> 2241.8823
$ws.Cells.Item(5, 13).Value = -2129.8823  # M5: -2214.8 -> -2129.8823

$ws.Cells.Item(120, 8).Value = 15232.5  # H120: 16837.143 -> 15232.5
$ws.Cells.Item(120, 9).Value = 10465  # I120: 12620 -> 10465
$ws.Cells.Item(120, 11).Value = 31395  # K120: 37860 -> 31395
$ws.Cells.Item(120, 13).Value = -26557  # M120: -33022 -> -26557

$ws.Cells.Item(122, 8).Value = 895.34784  # H122: 901.26086 -> 895.34784
$ws.Cells.Item(122, 9).Value = 420  # I122: 430 -> 420
$ws.Cells.Item(122, 10).Value = 1063.1177  # J122: 1067.5883 -> 1063.1177
$ws.Cells.Item(122, 11).Value = 3780  # K122: 3870 -> 3780
$ws.Cells.Item(122, 12).Value = 9568.059300000001  # L122: 9608.294699999999 -> 9568.059300000001
$ws.Cells.Item(122, 13).Value = -1330  # M122: -1420 -> -1330
$ws.Cells.Item(122, 14).Value = -14468.0593  # N122: -14508.2947 -> -14468.0593

$ws.Cells.Item(134, 8).Value = 1882.4  # H134: 2174.25 -> 1882.4
$ws.Cells.Item(134, 9).Value = 837  # I134: 1093.2 -> 837
$ws.Cells.Item(134, 10).Value = 3973.2  # J134: 3976 -> 3973.2
$ws.Cells.Item(134, 11).Value = 2511  # K134: 3279.6 -> 2511
$ws.Cells.Item(134, 12).Value = 11919.6  # L134: 11928 -> 11919.6
$ws.Cells.Item(134, 13).Value = 2559  # M134: 1790.4 -> 2559
$ws.Cells.Item(134, 14).Value = -22059.6  # N134: -22068 -> -22059.6

$ws.Cells.Item(135, 8).Value = 1489.5217  # H135: 1580.4286 -> 1489.5217
$ws.Cells.Item(135, 9).Value = 747.2941  # I135: 775.6 -> 747.2941
$ws.Cells.Item(135, 11).Value = 6725.6469  # K135: 6980.400000000001 -> 6725.6469
$ws.Cells.Item(135, 13).Value = -4190.6469  # M135: -4445.400000000001 -> -4190.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 50004.875  # H29: 8000 -> 50004.875
$ws.Cells.Item(29, 9).Value = 10003.5  # I29: 0 -> 10003.5
$ws.Cells.Item(29, 10).Value = 63338.668  # J29: 8000 -> 63338.668
$ws.Cells.Item(29, 11).Value = 10003.5  # K29: 0 -> 10003.5
$ws.Cells.Item(29, 12).Value = 63338.668  # L29: 8000 -> 63338.668
$ws.Cells.Item(29, 13).Value = -9713.5  # M29: None -> -9713.5
$ws.Cells.Item(29, 14).Value = -63918.668  # N29: -8580 -> -63918.668

$ws.Cells.Item(119, 8).Value = 38249.75  # H119: 37999.75 -> 38249.75
$ws.Cells.Item(119, 10).Value = 38249.75  # J119: 37999.75 -> 38249.75
$ws.Cells.Item(119, 12).Value = 38249.75  # L119: 37999.75 -> 38249.75
$ws.Cells.Item(119, 14).Value = -47925.75  # N119: -47675.75 -> -47925.75

$ws.Cells.Item(122, 8).Value = 3353.125  # H122: 3736.0476 -> 3353.125
$ws.Cells.Item(122, 9).Value = 2960  # I122: 4600.875 -> 2960
$ws.Cells.Item(122, 10).Value = 3531.818  # J122: 3203.8462 -> 3531.818
$ws.Cells.Item(122, 11).Value = 8880  # K122: 13802.625 -> 8880
$ws.Cells.Item(122, 12).Value = 10595.454  # L122: 9611.5386 -> 10595.454
$ws.Cells.Item(122, 13).Value = -6430  # M122: -11352.625 -> -6430
$ws.Cells.Item(122, 14).Value = -15495.454  # N122: -14511.5386 -> -15495.454

$ws.Cells.Item(132, 8).Value = 3126.842  # H132: 3293.5 -> 3126.842
$ws.Cells.Item(132, 9).Value = 1944.2222  # I132: 1966.4445 -> 1944.2222
$ws.Cells.Item(132, 10).Value = 4191.2  # J132: 4999.7144 -> 4191.2
$ws.Cells.Item(132, 11).Value = 5832.6666  # K132: 5899.333500000001 -> 5832.6666
$ws.Cells.Item(132, 12).Value = 12573.6  # L132: 14999.1432 -> 12573.6
$ws.Cells.Item(132, 13).Value = -3302.6666  # M132: -3369.333500000001 -> -3302.6666
$ws.Cells.Item(132, 14).Value = -17633.6  # N132: -20059.1432 -> -17633.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2414.92  # H82: 2393.074 -> 2414.92
$ws.Cells.Item(82, 9).Value = 1921.5385  # I82: 1945.7142 -> 1921.5385
$ws.Cells.Item(82, 10).Value = 2949.4167  # J82: 2874.8462 -> 2949.4167
$ws.Cells.Item(82, 11).Value = 1921.5385  # K82: 1945.7142 -> 1921.5385
$ws.Cells.Item(82, 12).Value = 2949.4167  # L82: 2874.8462 -> 2949.4167
$ws.Cells.Item(82, 13).Value = -1560.5385  # M82: -1584.7142 -> -1560.5385
$ws.Cells.Item(82, 14).Value = -3671.4167  # N82: -3596.8462 -> -3671.4167

$ws.Cells.Item(85, 8).Value = 2414.92  # H85: 2393.074 -> 2414.92
$ws.Cells.Item(85, 9).Value = 1921.5385  # I85: 1945.7142 -> 1921.5385
$ws.Cells.Item(85, 10).Value = 2949.4167  # J85: 2874.8462 -> 2949.4167
$ws.Cells.Item(85, 11).Value = 1921.5385  # K85: 1945.7142 -> 1921.5385
$ws.Cells.Item(85, 12).Value = 2949.4167  # L85: 2874.8462 -> 2949.4167
$ws.Cells.Item(85, 13).Value = -673.5385000000001  # M85: -697.7141999999999 -> -673.5385000000001
$ws.Cells.Item(85, 14).Value = -5445.4167  # N85: -5370.8462 -> -5445.4167

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 47346  # H31: 70019 -> 47346
$ws.Cells.Item(31, 9).Value = 2000  # I31: 0 -> 2000
$ws.Cells.Item(31, 11).Value = 2000  # K31: 0 -> 2000
$ws.Cells.Item(31, 13).Value = -1652  # M31: None -> -1652
